# Update Sheets via scheduled runner: refresh market price data cells
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")  # row 40
$ws.Range("H40").Value = 2363.853
$ws.Range("J40").Value = 2426.5715
$ws.Range("L40").Value = 2426.5715
$ws.Range("N40").Value = -2776.5715

$ws = $wb.Worksheets.Item("ALC")  # row 137
$ws.Range("H137").Value = 14723.412
$ws.Range("I137").Value = 8568.916999999999
$ws.Range("J137").Value = 18080.408
$ws.Range("K137").Value = 25706.751
$ws.Range("L137").Value = 54241.224
$ws.Range("M137").Value = -23156.751
$ws.Range("N137").Value = -59341.224

$ws = $wb.Worksheets.Item("ALC")  # row 138
$ws.Range("H138").Value = 5572.826
$ws.Range("J138").Value = 5470.769
$ws.Range("L138").Value = 16412.307
$ws.Range("N138").Value = -26692.307

$ws = $wb.Worksheets.Item("ARM")  # row 15
$ws.Range("H15").Value = 7499.5
$ws.Range("I15").Value = 1999
$ws.Range("J15").Value = 13000
$ws.Range("K15").Value = 1999
$ws.Range("L15").Value = 13000
$ws.Range("M15").Value = -1649
$ws.Range("N15").Value = -13700

$ws = $wb.Worksheets.Item("ARM")  # row 32
$ws.Range("H32").Value = 11876.671
$ws.Range("J32").Value = 33332.668
$ws.Range("L32").Value = 33332.668
$ws.Range("N32").Value = -33906.668

$ws = $wb.Worksheets.Item("ARM")  # row 61
$ws.Range("H61").Value = 12780.643
$ws.Range("I61").Value = 6179
$ws.Range("J61").Value = 21582.834
$ws.Range("K61").Value = 6179
$ws.Range("L61").Value = 21582.834
$ws.Range("M61").Value = -5967
$ws.Range("N61").Value = -22006.834

$ws = $wb.Worksheets.Item("ARM")  # row 74
$ws.Range("H74").Value = 8686.164000000001
$ws.Range("I74").Value = 7833.74
$ws.Range("K74").Value = 7833.74
$ws.Range("M74").Value = -6959.74

$ws = $wb.Worksheets.Item("ARM")  # row 77
$ws.Range("H77").Value = 8686.164000000001
$ws.Range("I77").Value = 7833.74
$ws.Range("K77").Value = 39168.7
$ws.Range("M77").Value = -34800.7

$ws = $wb.Worksheets.Item("ARM")  # row 110
$ws.Range("H110").Value = 1863.0476
$ws.Range("I110").Value = 1706.6316
$ws.Range("K110").Value = 1706.6316
$ws.Range("M110").Value = 338.3684000000001

$ws = $wb.Worksheets.Item("ARM")  # row 132
$ws.Range("H132").Value = 2567.3774
$ws.Range("I132").Value = 2533.5227
$ws.Range("K132").Value = 7600.5681
$ws.Range("M132").Value = -5070.5681

$ws = $wb.Worksheets.Item("ARM")  # row 136
$ws.Range("H136").Value = 12780.643
$ws.Range("I136").Value = 6179
$ws.Range("J136").Value = 21582.834
$ws.Range("K136").Value = 18537
$ws.Range("L136").Value = 64748.50199999999
$ws.Range("M136").Value = -15987
$ws.Range("N136").Value = -69848.50199999999

$ws = $wb.Worksheets.Item("BSM")  # row 19
$ws.Range("H19").Value = 6044
$ws.Range("J19").Value = 6044
$ws.Range("L19").Value = 6044
$ws.Range("N19").Value = -6390

$ws = $wb.Worksheets.Item("BSM")  # row 22
$ws.Range("H22").Value = 423.9
$ws.Range("I22").Value = 421.14285
$ws.Range("K22").Value = 421.14285
$ws.Range("M22").Value = -248.14285

$ws = $wb.Worksheets.Item("BSM")  # row 86
$ws.Range("H86").Value = 387818.66
$ws.Range("I86").Value = 770893.25
$ws.Range("K86").Value = 770893.25
$ws.Range("M86").Value = -769770.25

$ws = $wb.Worksheets.Item("BSM")  # row 89
$ws.Range("H89").Value = 387818.66
$ws.Range("I89").Value = 770893.25
$ws.Range("K89").Value = 3854466.25
$ws.Range("M89").Value = -3848850.25

$ws = $wb.Worksheets.Item("BSM")  # row 134
$ws.Range("H134").Value = 11616.412
$ws.Range("I134").Value = 5932.619
$ws.Range("K134").Value = 17797.857
$ws.Range("M134").Value = -15262.857

$ws = $wb.Worksheets.Item("CRP")  # row 7
$ws.Range("H7").Value = 36.833332
$ws.Range("I7").Value = 36.833332
$ws.Range("K7").Value = 36.833332
$ws.Range("M7").Value = 76.166668

$ws = $wb.Worksheets.Item("CRP")  # row 31
$ws.Range("H31").Value = 4628.5293
$ws.Range("I31").Value = 3741.7144
$ws.Range("K31").Value = 3741.7144
$ws.Range("M31").Value = -3446.7144

$ws = $wb.Worksheets.Item("CRP")  # row 34
$ws.Range("H34").Value = 4628.5293
$ws.Range("I34").Value = 3741.7144
$ws.Range("K34").Value = 3741.7144
$ws.Range("M34").Value = -3539.7144

$ws = $wb.Worksheets.Item("CRP")  # row 58
$ws.Range("H58").Value = 5556.0815
$ws.Range("I58").Value = 4460.72
$ws.Range("J58").Value = 6697.0835
$ws.Range("K58").Value = 4460.72
$ws.Range("L58").Value = 6697.0835
$ws.Range("M58").Value = -4257.72
$ws.Range("N58").Value = -7103.0835

$ws = $wb.Worksheets.Item("CRP")  # row 105
$ws.Range("H105").Value = 3308.5454
$ws.Range("I105").Value = 3439.4
$ws.Range("K105").Value = 3439.4
$ws.Range("M105").Value = -1692.4

$ws = $wb.Worksheets.Item("CRP")  # row 132
$ws.Range("H132").Value = 25973.475
$ws.Range("I132").Value = 19259.656
$ws.Range("K132").Value = 57778.96799999999
$ws.Range("M132").Value = -55248.96799999999

$ws = $wb.Worksheets.Item("CRP")  # row 134
$ws.Range("H134").Value = 3959.6553
$ws.Range("I134").Value = 2318.6924
$ws.Range("J134").Value = 5292.9375
$ws.Range("K134").Value = 6956.0772
$ws.Range("L134").Value = 15878.8125
$ws.Range("M134").Value = -4421.0772
$ws.Range("N134").Value = -20948.8125

$ws = $wb.Worksheets.Item("CRP")  # row 136
$ws.Range("H136").Value = 5556.0815
$ws.Range("I136").Value = 4460.72
$ws.Range("J136").Value = 6697.0835
$ws.Range("K136").Value = 13382.16
$ws.Range("L136").Value = 20091.2505
$ws.Range("M136").Value = -10832.16
$ws.Range("N136").Value = -25191.2505

$ws = $wb.Worksheets.Item("CUL")  # row 113
$ws.Range("H113").Value = 2202.55
$ws.Range("I113").Value = 2139.5557
$ws.Range("K113").Value = 6418.6671
$ws.Range("M113").Value = -4248.6671

$ws = $wb.Worksheets.Item("GSM")  # row 80
$ws.Range("H80").Value = 2026.375
$ws.Range("I80").Value = 2012.1666
$ws.Range("K80").Value = 2012.1666
$ws.Range("M80").Value = -1014.1666

$ws = $wb.Worksheets.Item("GSM")  # row 83
$ws.Range("H83").Value = 2026.375
$ws.Range("I83").Value = 2012.1666
$ws.Range("K83").Value = 10060.833
$ws.Range("M83").Value = -5068.833000000001

$ws = $wb.Worksheets.Item("GSM")  # row 120
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws = $wb.Worksheets.Item("LTW")  # row 13
$ws.Range("H13").Value = 423
$ws.Range("I13").Value = 6
$ws.Range("J13").Value = 840
$ws.Range("K13").Value = 6
$ws.Range("L13").Value = 840
$ws.Range("M13").Value = 134
$ws.Range("N13").Value = -1120

$ws = $wb.Worksheets.Item("LTW")  # row 40
$ws.Range("H40").Value = 9156.25
$ws.Range("I40").Value = 9464.571
$ws.Range("K40").Value = 9464.571
$ws.Range("M40").Value = -9328.571

$ws = $wb.Worksheets.Item("LTW")  # row 61
$ws.Range("H61").Value = 2712.1
$ws.Range("I61").Value = 2231.9285
$ws.Range("J61").Value = 3832.5
$ws.Range("K61").Value = 2231.9285
$ws.Range("L61").Value = 3832.5
$ws.Range("M61").Value = -2029.9285
$ws.Range("N61").Value = -4236.5

$ws = $wb.Worksheets.Item("LTW")  # row 113
$ws.Range("H113").Value = 2712.1
$ws.Range("I113").Value = 2231.9285
$ws.Range("J113").Value = 3832.5
$ws.Range("K113").Value = 2231.9285
$ws.Range("L113").Value = 3832.5
$ws.Range("M113").Value = -61.92849999999999
$ws.Range("N113").Value = -8172.5

$ws = $wb.Worksheets.Item("LTW")  # row 136
$ws.Range("H136").Value = 6183.947
$ws.Range("I136").Value = 5777.8354
$ws.Range("K136").Value = 17333.5062
$ws.Range("M136").Value = -14783.5062

$ws = $wb.Worksheets.Item("WVR")  # row 96
$ws.Range("H96").Value = 2210.923
$ws.Range("I96").Value = 1895.1666
$ws.Range("J96").Value = 6000
$ws.Range("K96").Value = 1895.1666
$ws.Range("L96").Value = 6000
$ws.Range("M96").Value = -522.1666
$ws.Range("N96").Value = -8746

$ws = $wb.Worksheets.Item("WVR")  # row 113
$ws.Range("H113").Value = 688.92
$ws.Range("I113").Value = 652.05554
$ws.Range("K113").Value = 1956.16662
$ws.Range("M113").Value = 213.83338

$ws = $wb.Worksheets.Item("WVR")  # row 122
$ws.Range("H122").Value = 3773
$ws.Range("I122").Value = 2754.5
$ws.Range("K122").Value = 8263.5
$ws.Range("M122").Value = -5813.5

$ws = $wb.Worksheets.Item("WVR")  # row 132
$ws.Range("H132").Value = 143225.58
$ws.Range("I132").Value = 189857.05
$ws.Range("K132").Value = 569571.1499999999
$ws.Range("M132").Value = -567041.1499999999

$ws = $wb.Worksheets.Item("WVR")  # row 136
$ws.Range("H136").Value = 6252321
$ws.Range("I136").Value = 10346307
$ws.Range("J136").Value = 3605.7896
$ws.Range("K136").Value = 31038921
$ws.Range("L136").Value = 10817.3688
$ws.Range("M136").Value = -31036371
$ws.Range("N136").Value = -15917.3688

